# Model accommodating both reading from files and from nasa:
# the spreadsheet no longer stores a pre-computed "m (kg)" mass column -
# that data/header is removed (now presumably computed/sourced elsewhere,
# e.g. from NASA data, at runtime rather than baked into the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "m (kg)" header text and the mass formulas/values in column H
# (rows 1 through 12), leaving the cell formatting/styles intact.
$ws.Range("H1:H12").ClearContents()

# Reflect the user's last selection: the whole of column H.
$ws.Range("H1:H1048576").Select()
